$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.343.30"
$ws.Range("E2").Value = "  +8.34%  "
$ws.Range("D3").Value = "2.617.88"
$ws.Range("E3").Value = "  +7.98%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'186.95"
$ws.Range("E5").Value = "  +16.00%  "
$ws.Range("D6").Value = "'587.93"
$ws.Range("E6").Value = "  +5.25%  "
$ws.Range("D8").Value = "'0.539"
$ws.Range("E8").Value = "  +5.19%  "
$ws.Range("E9").Value = "  +23.01%  "
$ws.Range("D10").Value = "2.618.95"
$ws.Range("E10").Value = "  +8.12%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "'0.364"
$ws.Range("E12").Value = "  +10.28%  "
$ws.Range("E13").Value = "  +4.35%  "
$ws.Range("E14").Value = "  +10.05%  "
$ws.Range("D15").Value = "74.243.64"
$ws.Range("E15").Value = "  +8.33%  "
$ws.Range("D16").Value = "3.087.29"
$ws.Range("E16").Value = "  +7.49%  "
$ws.Range("D17").Value = "'26.41"
$ws.Range("E17").Value = "  +14.49%  "
$ws.Range("D18").Value = "2.627.86"
$ws.Range("E18").Value = "  +8.14%  "
$ws.Range("D19").Value = "'9.19"
$ws.Range("E19").Value = "  +33.16%  "
$ws.Range("D20").Value = "'11.88"
$ws.Range("E20").Value = "  +13.66%  "
$ws.Range("D21").Value = "'375.28"
$ws.Range("E21").Value = "  +11.51%  "
$ws.Range("D22").Value = "'2.28"
$ws.Range("E22").Value = "  +18.95%  "
$ws.Range("E23").Value = "  +7.75%  "
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "'70.28"
$ws.Range("E25").Value = "  +5.24%  "
$ws.Range("D26").Value = "'4.21"
$ws.Range("E26").Value = "  +14.59%  "
$ws.Range("D27").Value = "'9.38"
$ws.Range("E27").Value = "  +14.68%  "
$ws.Range("D28").Value = "2.749.88"
$ws.Range("E28").Value = "  +7.57%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").Value = "0.0₃0955"
$ws.Range("E30").Value = "  +17.54%  "
$ws.Range("D31").Value = "'1.40"
$ws.Range("E31").Value = "  +22.60%  "
$ws.Range("D32").Value = "'8.01"
$ws.Range("E32").Value = "  +12.94%  "
$ws.Range("D33").Value = "'512.02"
$ws.Range("E33").Value = "  +20.25%  "
$ws.Range("E34").Value = "  +9.64%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").Value = "'0.123"
$ws.Range("E36").Value = "  +16.08%  "
$ws.Range("D37").Value = "'159.93"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "'19.26"
$ws.Range("E38").Value = "  +7.65%  "
$ws.Range("D39").Value = "'19.35"
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.0946"
$ws.Range("E41").Value = "  +32.33%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'4.95"
$ws.Range("E42").Value = "  +14.58%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.70"
$ws.Range("E43").Value = "  +13.61%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "'0.327"
$ws.Range("E44").Value = "  +10.24%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.42"
$ws.Range("E45").Value = "  +19.27%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'158.70"
$ws.Range("E46").Value = "  +21.32%  "
$ws.Range("D47").Value = "'1.19"
$ws.Range("E47").Value = "  +11.55%  "
$ws.Range("D48").Value = "'38.91"
$ws.Range("E48").Value = "  +4.02%  "
$ws.Range("D49").Value = "'3.65"
$ws.Range("E49").Value = "  +9.48%  "
$ws.Range("D50").Value = "'0.528"
$ws.Range("E50").Value = "  +10.35%  "
$ws.Range("D51").Value = "'20.51"
$ws.Range("E51").Value = "  +22.17%  "
